$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.644.61'
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").Value = '1.591.64'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").Value = '''210.30'
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D6").Value = '''0.515'
$ws.Range("E6").Value = '  +1.14%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("E8").Value = '  -0.43%  '

$ws.Range("E9").Value = '  -1.91%  '

$ws.Range("D10").Value = '''19.40'
$ws.Range("E10").Value = '  -1.30%  '

$ws.Range("D11").Value = '''0.0839'
$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").Value = '1.815.98'
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").Value = '1.597.54'
$ws.Range("E13").Value = '  -0.54%  '

$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("D15").Value = '''0.519'
$ws.Range("E15").Value = '  -1.60%  '

$ws.Range("D16").Value = '''64.29'
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").Value = '26.623.90'
$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").Value = '0.0₃0726'
$ws.Range("E18").Value = '  -0.26%  '

$ws.Range("E19").Value = '  +0.25%  '

$ws.Range("D20").Value = '''206.65'
$ws.Range("E20").Value = '  -0.56%  '

$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").Value = '''4.22'
$ws.Range("E22").Value = '  -0.64%  '

$ws.Range("D23").Value = '''2.30'
$ws.Range("E23").Value = '  -2.77%  '

$ws.Range("D25").Value = '''145.53'
$ws.Range("E25").Value = '  -1.14%  '

$ws.Range("E26").Value = '  +0.28%  '

$ws.Range("E27").Value = '  -2.35%  '

$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("D29").Value = '''15.20'
$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("D30").Value = '''0.0502'
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").Value = '''3.21'
$ws.Range("E32").Value = '  -1.32%  '

$ws.Range("E33").Value = '  +0.16%  '

$ws.Range("D34").Value = '''2.91'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").Value = '1.278.56'
$ws.Range("E35").Value = '  -3.80%  '

$ws.Range("E36").Value = '  +2.04%  '

$ws.Range("E37").Value = '  -1.74%  '

$ws.Range("E38").Value = '  -0.50%  '

$ws.Range("D39").Value = '''0.836'
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("E40").Value = '  +0.19%  '

$ws.Range("D41").Value = '''5.40'
$ws.Range("E41").Value = '  +0.42%  '

$ws.Range("E42").Value = '  +0.87%  '

$ws.Range("E43").Value = '  -0.01%  '

$ws.Range("D44").Value = '''63.11'
$ws.Range("E44").Value = '  -0.76%  '

$ws.Range("D45").Value = '1.728.55'
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("D46").Value = '''0.913'
$ws.Range("E46").Value = '  +9.34%  '

$ws.Range("D47").Value = '''89.72'
$ws.Range("E47").Value = '  -0.31%  '

$ws.Range("D48").Value = '''1.59'
$ws.Range("E48").Value = '  -0.91%  '

$ws.Range("E49").Value = '  +3.07%  '

$ws.Range("E50").Value = '  -0.90%  '

$ws.Range("E51").Value = '  +0.10%  '
